$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.452.30'
$ws.Range("E2").Value = '  -1.17%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.021.18'
$ws.Range("E3").Value = '  -0.96%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '585.84'
$ws.Range("E5").Value = '  -0.52%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.72'
$ws.Range("E6").Value = '  -3.24%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  -1.77%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.017.63'
$ws.Range("E9").Value = '  -1.25%  '
$ws.Range("E10").Value = '  -3.52%  '
$ws.Range("E11").Value = '  -0.57%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.450'
$ws.Range("E12").Value = '  +0.19%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000229'
$ws.Range("E13").Value = '  -2.77%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.70'
$ws.Range("E14").Value = '  -4.63%  '
$ws.Range("E15").Value = '  +2.12%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.518.96'
$ws.Range("E16").Value = '  -1.24%  '
$ws.Range("E17").Value = '  -1.16%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '62.417.46'
$ws.Range("E18").Value = '  -1.28%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.021.48'
$ws.Range("E19").Value = '  -1.18%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '462.91'
$ws.Range("E20").Value = '  -3.95%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.99'
$ws.Range("E21").Value = '  -3.16%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.687'
$ws.Range("E22").Value = '  -2.38%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.48'
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.64'
$ws.Range("E24").Value = '  -0.16%  '
$ws.Range("E25").Value = '  -4.69%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.36'
$ws.Range("E26").Value = '  -3.23%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.21'
$ws.Range("E27").Value = '  -2.42%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  +0.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("E30").Value = '  -1.51%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.12'
$ws.Range("E31").Value = '  -4.14%  '
$ws.Range("E32").Value = '  -4.92%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '28.49'
$ws.Range("E33").Value = '  +4.80%  '
$ws.Range("E34").Value = '  -1.49%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0₃0805'
$ws.Range("E35").Value = '  -1.46%  '
$ws.Range("E36").Value = '  -1.91%  '
$ws.Range("E37").Value = '  -3.56%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.12'
$ws.Range("E38").Value = '  -3.69%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '50.42'
$ws.Range("E39").Value = '  -0.09%  '
$ws.Range("E40").Value = '  -1.98%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.94'
$ws.Range("E41").Value = '  -8.79%  '
$ws.Range("E42").Value = '  +0.60%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '393.16'
$ws.Range("E43").Value = '  -10.03%  '
$ws.Range("E44").Value = '  -3.08%  '
$ws.Range("E45").Value = '  -0.76%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.758.66'
$ws.Range("E46").Value = '  -1.98%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '37.19'
$ws.Range("E47").Value = '  -5.32%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '128.83'
$ws.Range("E48").Value = '  -3.02%  '
$ws.Range("E50").Value = '  -0.23%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.21'
$ws.Range("E51").Value = '  -4.05%  '
